# Update NATMI ligand/receptor (Nodal-Acvr2b) TPM-derived metrics with
# refreshed values produced by the updated TPM pipeline.
# Columns: E=Ligand-expressing cells, F=Ligand detection rate,
# G=Ligand avg expr, H=Ligand total expr, I/J=Ligand derived specificity,
# M=Receptor avg expr, N=Receptor total expr, O/P=Receptor derived specificity,
# Q=Edge avg expr weight, R=Edge total expr weight,
# S/T=Edge derived specificity (avg/total).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Sending cluster ECs -> Target cluster ECs
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 1.152264333333333
$ws.Cells.Item(2, 8).Value = 3.456793
$ws.Cells.Item(2, 9).Value = 0.3930660006090215
$ws.Cells.Item(2, 10).Value = 0.3930660006090216
$ws.Cells.Item(2, 13).Value = 1.485259333333333
$ws.Cells.Item(2, 14).Value = 4.455778
$ws.Cells.Item(2, 15).Value = 0.3057455162066235
$ws.Cells.Item(2, 16).Value = 0.3057455162066235
$ws.Cells.Item(2, 17).Value = 1.711411355550445
$ws.Cells.Item(2, 18).Value = 15.402702199954
$ws.Cells.Item(2, 19).Value = 0.1201781672594783
$ws.Cells.Item(2, 20).Value = 0.1201781672594783

# Row 3: Sending cluster ECs -> Target cluster FAPs
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 1.152264333333333
$ws.Cells.Item(3, 8).Value = 3.456793
$ws.Cells.Item(3, 9).Value = 0.3930660006090215
$ws.Cells.Item(3, 10).Value = 0.3930660006090216
$ws.Cells.Item(3, 15).Value = 0.2805555239151429
$ws.Cells.Item(3, 16).Value = 0.2805555239151429
$ws.Cells.Item(3, 17).Value = 1.570410305432889
$ws.Cells.Item(3, 18).Value = 14.133692748896
$ws.Cells.Item(3, 19).Value = 0.1102768377340939
$ws.Cells.Item(3, 20).Value = 0.1102768377340939

# Row 4: Sending cluster ECs -> Target cluster MuSCs
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 1.152264333333333
$ws.Cells.Item(4, 8).Value = 3.456793
$ws.Cells.Item(4, 9).Value = 0.3930660006090215
$ws.Cells.Item(4, 10).Value = 0.3930660006090216
$ws.Cells.Item(4, 15).Value = 0.4136989598782336
$ws.Cells.Item(4, 16).Value = 0.4136989598782336
$ws.Cells.Item(4, 17).Value = 2.315681049060889
$ws.Cells.Item(4, 18).Value = 20.841129441548
$ws.Cells.Item(4, 19).Value = 0.1626109956154493
$ws.Cells.Item(4, 20).Value = 0.1626109956154494

# Row 5: Sending cluster FAPs -> Target cluster ECs
$ws.Cells.Item(5, 9).Value = 0.4761983545501621
$ws.Cells.Item(5, 10).Value = 0.476198354550162
$ws.Cells.Item(5, 13).Value = 1.485259333333333
$ws.Cells.Item(5, 14).Value = 4.455778
$ws.Cells.Item(5, 15).Value = 0.3057455162066235
$ws.Cells.Item(5, 16).Value = 0.3057455162066235
$ws.Cells.Item(5, 17).Value = 2.073370045256667
$ws.Cells.Item(5, 18).Value = 18.66033040731
$ws.Cells.Item(5, 19).Value = 0.145595511728684
$ws.Cells.Item(5, 20).Value = 0.145595511728684

# Row 6: Sending cluster FAPs -> Target cluster FAPs
$ws.Cells.Item(6, 9).Value = 0.4761983545501621
$ws.Cells.Item(6, 10).Value = 0.476198354550162
$ws.Cells.Item(6, 15).Value = 0.2805555239151429
$ws.Cells.Item(6, 16).Value = 0.2805555239151429
$ws.Cells.Item(6, 19).Value = 0.1336000788483497
$ws.Cells.Item(6, 20).Value = 0.1336000788483497

# Row 7: Sending cluster FAPs -> Target cluster MuSCs
$ws.Cells.Item(7, 9).Value = 0.4761983545501621
$ws.Cells.Item(7, 10).Value = 0.476198354550162
$ws.Cells.Item(7, 15).Value = 0.4136989598782336
$ws.Cells.Item(7, 16).Value = 0.4136989598782336
$ws.Cells.Item(7, 19).Value = 0.1970027639731284
$ws.Cells.Item(7, 20).Value = 0.1970027639731284

# Row 8: Sending cluster MuSCs -> Target cluster ECs
$ws.Cells.Item(8, 9).Value = 0.1307356448408163
$ws.Cells.Item(8, 10).Value = 0.1307356448408163
$ws.Cells.Item(8, 13).Value = 1.485259333333333
$ws.Cells.Item(8, 14).Value = 4.455778
$ws.Cells.Item(8, 15).Value = 0.3057455162066235
$ws.Cells.Item(8, 16).Value = 0.3057455162066235
$ws.Cells.Item(8, 17).Value = 0.5692236591542222
$ws.Cells.Item(8, 18).Value = 5.123012932388
$ws.Cells.Item(8, 19).Value = 0.03997183721846119
$ws.Cells.Item(8, 20).Value = 0.03997183721846119

# Row 9: Sending cluster MuSCs -> Target cluster FAPs
$ws.Cells.Item(9, 9).Value = 0.1307356448408163
$ws.Cells.Item(9, 10).Value = 0.1307356448408163
$ws.Cells.Item(9, 15).Value = 0.2805555239151429
$ws.Cells.Item(9, 16).Value = 0.2805555239151429
$ws.Cells.Item(9, 19).Value = 0.03667860733269928
$ws.Cells.Item(9, 20).Value = 0.03667860733269928

# Row 10: Sending cluster MuSCs -> Target cluster MuSCs
$ws.Cells.Item(10, 9).Value = 0.1307356448408163
$ws.Cells.Item(10, 10).Value = 0.1307356448408163
$ws.Cells.Item(10, 15).Value = 0.4136989598782336
$ws.Cells.Item(10, 16).Value = 0.4136989598782336
$ws.Cells.Item(10, 19).Value = 0.05408520028965588
$ws.Cells.Item(10, 20).Value = 0.05408520028965588
